$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update res_bus vm_pu values for Case_4_221 (380 kV case)

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037520936406501
$ws.Range("D2").Value = 1.038022481527427
$ws.Range("E2").Value = 1.050988746968421
$ws.Range("F2").Value = 1.058247659740094
$ws.Range("I2").Value = 1.034808401007407
$ws.Range("J2").Value = 1.042623357493824
$ws.Range("K2").Value = 1.04081169903153
$ws.Range("L2").Value = 1.053741495435823
$ws.Range("M2").Value = 1.06098041446597
$ws.Range("N2").Value = 1.044104002252036

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038579583786722
$ws.Range("D3").Value = 1.038774347677602
$ws.Range("E3").Value = 1.052031534457275
$ws.Range("F3").Value = 1.05935729717615
$ws.Range("I3").Value = 1.034986048911102
$ws.Range("J3").Value = 1.043325869883406
$ws.Range("K3").Value = 1.041373646636141
$ws.Range("L3").Value = 1.054596246794427
$ws.Range("M3").Value = 1.061903297771505
$ws.Range("N3").Value = 1.044807512289791

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.039264891991229
$ws.Range("D4").Value = 1.039260874552283
$ws.Range("E4").Value = 1.052706933106319
$ws.Range("F4").Value = 1.060075995463201
$ws.Range("I4").Value = 1.035099683966695
$ws.Range("J4").Value = 1.043780164854523
$ws.Range("K4").Value = 1.041736620706805
$ws.Range("L4").Value = 1.055149369150278
$ws.Range("M4").Value = 1.062500561331421
$ws.Range("N4").Value = 1.045262452411875

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039553065630867
$ws.Range("D5").Value = 1.039465414258577
$ws.Range("E5").Value = 1.052991024829326
$ws.Range("F5").Value = 1.06037830078217
$ws.Range("I5").Value = 1.03514714103058
$ws.Range("J5").Value = 1.043971083695989
$ws.Range("K5").Value = 1.041889060443896
$ws.Range("L5").Value = 1.055381911281308
$ws.Range("M5").Value = 1.062751673669821
$ws.Range("N5").Value = 1.045453642379994

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039601455362971
$ws.Range("D6").Value = 1.039499757575667
$ws.Range("E6").Value = 1.053038734106943
$ws.Range("F6").Value = 1.060429068816359
$ws.Range("I6").Value = 1.035155090793446
$ws.Range("J6").Value = 1.044003135885041
$ws.Range("K6").Value = 1.04191464666893
$ws.Range("L6").Value = 1.055420956675171
$ws.Range("M6").Value = 1.062793837860261
$ws.Range("N6").Value = 1.045485740086831

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.039268742306127
$ws.Range("D7").Value = 1.039263607607879
$ws.Range("E7").Value = 1.052710728548672
$ws.Range("F7").Value = 1.060080034235446
$ws.Range("I7").Value = 1.035100319328854
$ws.Range("J7").Value = 1.043782716184172
$ws.Range("K7").Value = 1.041738658221962
$ws.Range("L7").Value = 1.055152476352738
$ws.Range("M7").Value = 1.062503916617842
$ws.Range("N7").Value = 1.045265007364704

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037878650686493
$ws.Range("D8").Value = 1.038276573764273
$ws.Range("E8").Value = 1.051341028014773
$ws.Range("F8").Value = 1.058622524127461
$ws.Range("I8").Value = 1.034868710053213
$ws.Range("J8").Value = 1.042860832111682
$ws.Range("K8").Value = 1.041001744828907
$ws.Range("L8").Value = 1.054030353950618
$ws.Range("M8").Value = 1.061292287052138
$ws.Range("N8").Value = 1.04434181411109

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035431370558803
$ws.Range("D9").Value = 1.036537475596241
$ws.Range("E9").Value = 1.048932400348085
$ws.Range("F9").Value = 1.056059492770952
$ws.Range("I9").Value = 1.034450526428367
$ws.Range("J9").Value = 1.041234243209691
$ws.Range("K9").Value = 1.039698301055736
$ws.Range("L9").Value = 1.052053360176231
$ws.Range("M9").Value = 1.059157993074377
$ws.Range("N9").Value = 1.042712915266286

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.03380135047047
$ws.Range("D10").Value = 1.035378243899118
$ws.Range("E10").Value = 1.04733000308724
$ws.Range("F10").Value = 1.054354375627253
$ws.Range("I10").Value = 1.034164991376717
$ws.Range("J10").Value = 1.040148448168418
$ws.Range("K10").Value = 1.038826065406603
$ws.Range("L10").Value = 1.050735603379658
$ws.Range("M10").Value = 1.057735652920729
$ws.Range("N10").Value = 1.041625578271512

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033095887457372
$ws.Range("D11").Value = 1.034876332393335
$ws.Range("E11").Value = 1.046636946436279
$ws.Range("F11").Value = 1.053616890296122
$ws.Range("I11").Value = 1.034039754790254
$ws.Range("J11").Value = 1.039677957284567
$ws.Range("K11").Value = 1.038447606576169
$ws.Range("L11").Value = 1.050165059855922
$ws.Range("M11").Value = 1.05711989064174
$ws.Range("N11").Value = 1.04115441923664

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.032833898844159
$ws.Range("D12").Value = 1.034689907306745
$ws.Range("E12").Value = 1.046379633499904
$ws.Range("F12").Value = 1.053343081944
$ws.Range("I12").Value = 1.033992996481444
$ws.Range("J12").Value = 1.039503145934031
$ws.Range("K12").Value = 1.038306914025243
$ws.Range("L12").Value = 1.049953142835562
$ws.Range("M12").Value = 1.056891187331751
$ws.Range("N12").Value = 1.040979359633935

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032890093953347
$ws.Range("D13").Value = 1.03472989578446
$ws.Range("E13").Value = 1.046434822602938
$ws.Range("F13").Value = 1.053401809037501
$ws.Range("I13").Value = 1.034003037148226
$ws.Range("J13").Value = 1.03954064582609
$ws.Range("K13").Value = 1.038337098302719
$ws.Range("L13").Value = 1.049998599380058
$ws.Range("M13").Value = 1.056940244141624
$ws.Range("N13").Value = 1.041016912780142

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033074230342627
$ws.Range("D14").Value = 1.034860922286928
$ws.Range("E14").Value = 1.046615674436754
$ws.Range("F14").Value = 1.053594254630492
$ws.Range("I14").Value = 1.034035894625743
$ws.Range("J14").Value = 1.03966350834463
$ws.Range("K14").Value = 1.038435979252203
$ws.Range("L14").Value = 1.050147542559754
$ws.Range("M14").Value = 1.057100985582407
$ws.Range("N14").Value = 1.041139949777551

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.033187689699567
$ws.Range("D15").Value = 1.034941653012655
$ws.Range("E15").Value = 1.04672711901607
$ws.Range("F15").Value = 1.053712843481851
$ws.Range("I15").Value = 1.034056107424332
$ws.Range("J15").Value = 1.039739201365049
$ws.Range("K15").Value = 1.038496887696041
$ws.Range("L15").Value = 1.050239312466102
$ws.Range("M15").Value = 1.05720002610709
$ws.Range("N15").Value = 1.041215750290741

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.033848177001677
$ws.Range("D16").Value = 1.035411555085265
$ws.Range("E16").Value = 1.047376015675377
$ws.Range("F16").Value = 1.054403337839825
$ws.Range("I16").Value = 1.034173269266773
$ws.Range("J16").Value = 1.040179666065844
$ws.Range("K16").Value = 1.038851166173483
$ws.Range("L16").Value = 1.050773469622367
$ws.Range("M16").Value = 1.057776521582771
$ws.Range("N16").Value = 1.041656840501933

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034262575784046
$ws.Range("D17").Value = 1.035706324347639
$ws.Range("E17").Value = 1.047783263558617
$ws.Range("F17").Value = 1.054836692102062
$ws.Range("I17").Value = 1.034246334114068
$ws.Range("J17").Value = 1.040455868503853
$ws.Range("K17").Value = 1.039073188450043
$ws.Range("L17").Value = 1.051108546824596
$ws.Range("M17").Value = 1.058138174287472
$ws.Range("N17").Value = 1.041933435179084

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034504321028451
$ws.Range("D18").Value = 1.035878262420985
$ws.Range("E18").Value = 1.048020880756945
$ws.Range("F18").Value = 1.055089541511808
$ws.Range("I18").Value = 1.034288797402249
$ws.Range("J18").Value = 1.040616940430933
$ws.Range("K18").Value = 1.039202615484162
$ws.Range("L18").Value = 1.051303996874296
$ws.Range("M18").Value = 1.058349131798872
$ws.Range("N18").Value = 1.042094735846776

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034586755591338
$ws.Range("D19").Value = 1.035936889544366
$ws.Range("E19").Value = 1.048101915033912
$ws.Range("F19").Value = 1.055175770396637
$ws.Range("I19").Value = 1.034303250124289
$ws.Range("J19").Value = 1.040671856315089
$ws.Range("K19").Value = 1.039246734059477
$ws.Range("L19").Value = 1.051370641165417
$ws.Range("M19").Value = 1.058421064823542
$ws.Range("N19").Value = 1.042149729717787

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.03421811125904
$ws.Range("D20").Value = 1.035674697961373
$ws.Range("E20").Value = 1.047739561811265
$ws.Range("F20").Value = 1.054790188891261
$ws.Range("I20").Value = 1.034238510895551
$ws.Range("J20").Value = 1.04042623793404
$ws.Range("K20").Value = 1.039049375295646
$ws.Range("L20").Value = 1.051072595674239
$ws.Range("M20").Value = 1.058099371173259
$ws.Range("N20").Value = 1.041903762530464

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03302000533891
$ws.Range("D21").Value = 1.034822338043782
$ws.Range("E21").Value = 1.046562414788373
$ws.Range("F21").Value = 1.053537580706183
$ws.Range("I21").Value = 1.034026225537219
$ws.Range("J21").Value = 1.039627329776544
$ws.Range("K21").Value = 1.038406864478971
$ws.Range("L21").Value = 1.050103682276709
$ws.Range("M21").Value = 1.057053650730534
$ws.Range("N21").Value = 1.041103719831748

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.032267007925574
$ws.Range("D22").Value = 1.034286467484786
$ws.Range("E22").Value = 1.045822985261584
$ws.Range("F22").Value = 1.052750748371074
$ws.Range("I22").Value = 1.033891365635269
$ws.Range("J22").Value = 1.039124734798831
$ws.Range("K22").Value = 1.038002220937599
$ws.Range("L22").Value = 1.04949453576515
$ws.Range("M22").Value = 1.056396270248946
$ws.Range("N22").Value = 1.040600411111518

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.032666157823615
$ws.Range("D23").Value = 1.034570538409468
$ws.Range("E23").Value = 1.046214905423532
$ws.Range("F23").Value = 1.053167793649846
$ws.Range("I23").Value = 1.033962988882182
$ws.Range("J23").Value = 1.039391197312708
$ws.Range("K23").Value = 1.038216793667447
$ws.Range("L23").Value = 1.049817451273433
$ws.Range("M23").Value = 1.056744750039577
$ws.Range("N23").Value = 1.040867252032729

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034238202760407
$ws.Range("D24").Value = 1.035688988551166
$ws.Range("E24").Value = 1.047759308513616
$ws.Range("F24").Value = 1.054811201439493
$ws.Range("I24").Value = 1.034242046347694
$ws.Range("J24").Value = 1.040439626812772
$ws.Range("K24").Value = 1.039060135666185
$ws.Range("L24").Value = 1.051088840435465
$ws.Range("M24").Value = 1.058116904594338
$ws.Range("N24").Value = 1.041917170422939

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.036063786371205
$ws.Range("D25").Value = 1.036987047997426
$ws.Range("E25").Value = 1.049554498541121
$ws.Range("F25").Value = 1.056721469627842
$ws.Range("I25").Value = 1.03455982717472
$ws.Range("J25").Value = 1.041655003876937
$ws.Range("K25").Value = 1.040035851393672
$ws.Range("L25").Value = 1.052564419619012
$ws.Range("M25").Value = 1.059709668369984
$ws.Range("N25").Value = 1.043134273461941

Write-Host "Updated vm_pu values for Case_4_221 (380 kV) across rows 2-25"
